# Update BP terminal gate pricing sheet with the latest daily figures.
# The workbook rolls each state's price table forward by one day: new rows
# are stamped with the latest effective date and updated Diesel/ULP/PULP/e10
# prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New South Wales (Sydney-Botany / Sydney-Silverwater / Newcastle) ---
$ws.Range("A8").Value  = 45951
$ws.Range("D8").Value  = 160.47
$ws.Range("E8").Value  = 157.9
$ws.Range("F8").Value  = 167.9
$ws.Range("G8").Value  = 158.06

$ws.Range("A9").Value  = 45951
$ws.Range("D9").Value  = 160.47
$ws.Range("E9").Value  = 157.9
$ws.Range("F9").Value  = 167.9
$ws.Range("G9").Value  = 158.06

$ws.Range("A10").Value = 45951
$ws.Range("D10").Value = 162.69
$ws.Range("E10").Value = 160.26
$ws.Range("F10").Value = 170.26
$ws.Range("G10").Value = 160.73

$ws.Range("A11").Value = 45948
$ws.Range("D11").Value = 161.22
$ws.Range("E11").Value = 158.03
$ws.Range("F11").Value = 168.03
$ws.Range("G11").Value = 158.19

$ws.Range("A12").Value = 45948
$ws.Range("D12").Value = 161.22
$ws.Range("E12").Value = 158.03
$ws.Range("F12").Value = 168.03
$ws.Range("G12").Value = 158.19

$ws.Range("A13").Value = 45948
$ws.Range("D13").Value = 163.59
$ws.Range("E13").Value = 160.41
$ws.Range("F13").Value = 170.41
$ws.Range("G13").Value = 160.88

# --- Northern Territory (Darwin) ---
$ws.Range("A17").Value = 45951
$ws.Range("D17").Value = 166.12
$ws.Range("E17").Value = 163.09
$ws.Range("F17").Value = 173.09

$ws.Range("A18").Value = 45948
$ws.Range("D18").Value = 166.9
$ws.Range("E18").Value = 163.22
$ws.Range("F18").Value = 173.22

# --- Queensland (Brisbane / Cairns / Gladstone / Mackay / Townsville) ---
$ws.Range("A22").Value = 45951
$ws.Range("D22").Value = 161.39
$ws.Range("E22").Value = 159.15
$ws.Range("F22").Value = 168.75
$ws.Range("G22").Value = 160.32

$ws.Range("A23").Value = 45951
$ws.Range("D23").Value = 167.46
$ws.Range("E23").Value = 163.86
$ws.Range("F23").Value = 173.86

$ws.Range("A24").Value = 45951
$ws.Range("D24").Value = 167.26
$ws.Range("E24").Value = 164.07
$ws.Range("F24").Value = 174.07

$ws.Range("A25").Value = 45951
$ws.Range("D25").Value = 168.1
$ws.Range("E25").Value = 163.46
$ws.Range("F25").Value = 173.46
$ws.Range("G25").Value = 163.29

$ws.Range("A26").Value = 45951
$ws.Range("D26").Value = 166.81
$ws.Range("E26").Value = 165
$ws.Range("F26").Value = 175

$ws.Range("A27").Value = 45948
$ws.Range("D27").Value = 162.25
$ws.Range("E27").Value = 159.3
$ws.Range("F27").Value = 168.9
$ws.Range("G27").Value = 160.48

$ws.Range("A28").Value = 45948
$ws.Range("D28").Value = 168.25
$ws.Range("E28").Value = 164.01
$ws.Range("F28").Value = 174.01

$ws.Range("A29").Value = 45948
$ws.Range("D29").Value = 168.05
$ws.Range("E29").Value = 164.23
$ws.Range("F29").Value = 174.23

$ws.Range("A30").Value = 45948
$ws.Range("D30").Value = 168.88
$ws.Range("E30").Value = 163.62
$ws.Range("F30").Value = 173.62
$ws.Range("G30").Value = 163.45

$ws.Range("A31").Value = 45948
$ws.Range("D31").Value = 167.6
$ws.Range("E31").Value = 165.16
$ws.Range("F31").Value = 175.16

# --- South Australia (Adelaide) ---
$ws.Range("A35").Value = 45951
$ws.Range("D35").Value = 161.05
$ws.Range("E35").Value = 157.36
$ws.Range("F35").Value = 166.36

$ws.Range("A36").Value = 45948
$ws.Range("D36").Value = 161.95
$ws.Range("E36").Value = 157.52
$ws.Range("F36").Value = 166.52

# --- Tasmania (Burnie / Hobart) ---
$ws.Range("A40").Value = 45951
$ws.Range("D40").Value = 166.6
$ws.Range("E40").Value = 162.84
$ws.Range("F40").Value = 172.84

$ws.Range("A41").Value = 45951
$ws.Range("D41").Value = 166.31
$ws.Range("E41").Value = 163.26
$ws.Range("F41").Value = 173.26

$ws.Range("A42").Value = 45948
$ws.Range("D42").Value = 167.39
$ws.Range("E42").Value = 162.99
$ws.Range("F42").Value = 172.99

$ws.Range("A43").Value = 45948
$ws.Range("D43").Value = 167.1
$ws.Range("E43").Value = 163.41
$ws.Range("F43").Value = 173.41

# --- Victoria (Geelong / Melbourne) ---
$ws.Range("A47").Value = 45951
$ws.Range("D47").Value = 161.34
$ws.Range("E47").Value = 159.06
$ws.Range("F47").Value = 169.06

$ws.Range("A48").Value = 45951
$ws.Range("D48").Value = 161.32
$ws.Range("E48").Value = 159.23
$ws.Range("F48").Value = 169.23

$ws.Range("A49").Value = 45948
$ws.Range("D49").Value = 162.05
$ws.Range("E49").Value = 159.28
$ws.Range("F49").Value = 169.28

$ws.Range("A50").Value = 45948
$ws.Range("D50").Value = 162.03
$ws.Range("E50").Value = 159.45
$ws.Range("F50").Value = 169.45

# --- Western Australia (Broome / Esperance / Geraldton / Kalgoorlie / Kwinana) ---
$ws.Range("A54").Value = 45951
$ws.Range("D54").Value = 176.77
$ws.Range("E54").Value = 173.16
$ws.Range("F54").Value = 183.16

$ws.Range("A55").Value = 45951
$ws.Range("D55").Value = 164.43
$ws.Range("E55").Value = 170.52
$ws.Range("F55").Value = 180.52

$ws.Range("A56").Value = 45951
$ws.Range("D56").Value = 166.82

$ws.Range("A57").Value = 45951
$ws.Range("D57").Value = 166.48
$ws.Range("E57").Value = 164.79

$ws.Range("A58").Value = 45951
$ws.Range("D58").Value = 162.39
$ws.Range("E58").Value = 160.84
$ws.Range("F58").Value = 170.84

$ws.Range("A59").Value = 45951
$ws.Range("D59").Value = 169.1
$ws.Range("E59").Value = 171.33

$ws.Range("A60").Value = 45948
$ws.Range("D60").Value = 177.54
$ws.Range("E60").Value = 173.34
$ws.Range("F60").Value = 183.34

$ws.Range("A61").Value = 45948
$ws.Range("D61").Value = 165.22
$ws.Range("E61").Value = 170.78
$ws.Range("F61").Value = 180.78

$ws.Range("A62").Value = 45948
$ws.Range("D62").Value = 167.72

$ws.Range("A63").Value = 45948
$ws.Range("D63").Value = 167.37
$ws.Range("E63").Value = 165.05

$ws.Range("A64").Value = 45948
$ws.Range("D64").Value = 163.28
$ws.Range("E64").Value = 161.1
$ws.Range("F64").Value = 171.1

$ws.Range("A65").Value = 45948
$ws.Range("D65").Value = 169.88
$ws.Range("E65").Value = 171.5

$wb.Save()
